$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LanguageCandidates")

for ($row = 2; $row -le 24; $row++) {
    $predicateCell = $ws.Cells.Item($row, 9)   # Column I: PredictionPredicates
    $identityCell  = $ws.Cells.Item($row, 12)  # Column L: HasIdentity

    $currentText = $predicateCell.Value2
    $hasIdentity = $identityCell.Value2

    if ($hasIdentity -eq 1) {
        $suffix = ", Has Identity"
    } else {
        $suffix = ", Has no Identity"
    }

    $predicateCell.Value2 = "$currentText$suffix"
}
